# imb_perluasan_contoh.xlsx — "fix bug in imb and rekap"
#
# The header in column K was mislabeled "Fungsi Kegiatan" (Activity Function)
# when it should read "Fungsi Bangunan" (Building Function). Correct it, and
# leave the active selection on that cell (matching the author's last-saved
# cursor position).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the header text in K1.
$ws.Range("K1").Value = "Fungsi Bangunan"

# Match the saved selection/active cell.
$ws.Range("K1").Select()
